$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 77740542.77
$ws.Range("P2").Value = 170.707163872
$ws.Range("Q2").Value = 243931772.46
$ws.Range("R2").Value = 535.6394433482
$ws.Range("S2").Value = 56671647.01
$ws.Range("T2").Value = 124.4428683969
$ws.Range("U2").Value = -50460802.74
$ws.Range("V2").Value = -110.8047386283
$ws.Range("W2").Value = 109932.49
$ws.Range("X2").Value = 0.2413960968
$ws.Range("Y2").Value = 50600735.23
$ws.Range("Z2").Value = 111.1120104539
$ws.Range("AA2").Value = 18269038.64
$ws.Range("AB2").Value = 40.1162078599
$ws.Range("AC2").Value = 45540293.1
$ws.Range("AD2").Value = 21714.4438516938
